# Add additional volunteer information to the Volunteer Information sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 60 already had an entry (email "saitejavankadari2009@gmail.com" / school
# "Ascension of Our Lord") that was missing the volunteer's name and contact
# number. Fill those two cells in.
$ws.Range("A60").Value = "Sai"
$ws.Range("E60").Value = 6477798466

# Row 61 is a brand new volunteer record.
$ws.Range("A61").Value = "Janarth Kulenthiran"
$ws.Range("B61").Value = "kulenthirankk@hotmail.com"

# Match the formatting already used throughout column D (School) for the rest
# of the table by copying the format from the cell directly above before
# setting the new value.
$ws.Range("D60").Copy()
$ws.Range("D61").PasteSpecial(-4122)
$ws.Range("D61").Value = "Stephen Lewis Secondary"

$ws.Range("E61").Value = "905-813-9777 / 4167329912"

# Match the row height used by the rest of the table rows.
$ws.Rows.Item(61).RowHeight = 18.75

# Leave the selection on the last cell touched, like a user would after
# finishing data entry.
$ws.Range("E61").Select()
